# The source workbook was authored in an older Excel build and is being
# brought up to date: the sheet is rescrolled/reselected to where the
# author was last working (row ~15, cell D52) and the data column width
# is nudged from 16.57 to 16.5 characters, matching the new file's
# <col .../> width and <selection activeCell="D52" sqref="D52"/>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:M: 16.5703125 -> 16.5 (character width, stored in xl/worksheets/sheet1.xml)
$ws.Columns("A:M").ColumnWidth = 15.67

# Scroll the view down and move the selection to D52, like the new
# <sheetView topLeftCell="A15"><selection activeCell="D52" sqref="D52"/>
$excel.Goto($ws.Range("A15"), $true)
$ws.Range("D52").Select()
